$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.893.13"
$ws.Range("E2").Value = "  -0.12%  "
# Row 3
$ws.Range("D3").Value = "1.813.62"
$ws.Range("E3").Value = "  +1.53%  "
# Row 4
$ws.Range("E4").Value = "  -0.64%  "
# Row 5
$ws.Range("D5").Value = "'311.29"
$ws.Range("E5").Value = "  +0.68%  "
# Row 6
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  -0.66%  "
# Row 7
$ws.Range("D7").Value = "'0.4288"
$ws.Range("E7").Value = "  +1.30%  "
# Row 8
$ws.Range("D8").Value = "'0.3688"
$ws.Range("E8").Value = "  +2.32%  "
# Row 9
$ws.Range("E9").Value = "  +1.42%  "
# Row 10
$ws.Range("D10").Value = "'0.8615"
$ws.Range("E10").Value = "  +2.43%  "
# Row 11
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'21.05"
$ws.Range("E11").Value = "  +3.88%  "
# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "2.013.39"
$ws.Range("E12").Value = "  +12.08%  "
# Row 13
$ws.Range("E13").Value = "  +4.82%  "
# Row 14
$ws.Range("D14").Value = "'5.396"
$ws.Range("E14").Value = "  +2.88%  "
# Row 15
$ws.Range("D15").Value = "'0.06897"
$ws.Range("E15").Value = "  +0.79%  "
# Row 16
$ws.Range("D16").Value = "'80.66"
$ws.Range("E16").Value = "  +1.61%  "
# Row 17
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  -1.11%  "
# Row 18
$ws.Range("D18").Value = "'0.000008932"
$ws.Range("E18").Value = "  +3.03%  "
# Row 19
$ws.Range("D19").Value = "'1.004"
$ws.Range("E19").Value = "  -0.63%  "
# Row 20
$ws.Range("D20").Value = "'15.17"
$ws.Range("E20").Value = "  +1.82%  "
# Row 21
$ws.Range("D21").Value = "26.943.37"
$ws.Range("E21").Value = "  -0.03%  "
# Row 22
$ws.Range("D22").Value = "'5.178"
$ws.Range("E22").Value = "  +2.65%  "
# Row 23
$ws.Range("E23").Value = "  +0.27%  "
# Row 24
$ws.Range("D24").Value = "2.249.19"
$ws.Range("E24").Value = "  +10.94%  "
# Row 25
$ws.Range("D25").Value = "'153.59"
$ws.Range("E25").Value = "  +0.19%  "
# Row 26
$ws.Range("E26").Value = "  -3.15%  "
# Row 27
$ws.Range("D27").Value = "'18.22"
$ws.Range("E27").Value = "  +0.24%  "
# Row 28
$ws.Range("D28").Value = "'5.213"
$ws.Range("E28").Value = "  +4.13%  "
# Row 29
$ws.Range("D29").Value = "'114.89"
$ws.Range("E29").Value = "  +0.55%  "
# Row 30
$ws.Range("D30").Value = "'1.870"
$ws.Range("E30").Value = "  +15.35%  "
# Row 31
$ws.Range("D31").Value = "'0.08930"
$ws.Range("E31").Value = "  -0.08%  "
# Row 32
$ws.Range("D32").Value = "'0.7426"
$ws.Range("E32").Value = "  +3.34%  "
# Row 33
$ws.Range("D33").Value = "'1.163"
$ws.Range("E33").Value = "  +7.37%  "
# Row 34
$ws.Range("D34").Value = "'4.416"
$ws.Range("E34").Value = "  +2.55%  "
# Row 35
$ws.Range("D35").Value = "'2.797"
$ws.Range("E35").Value = "  -1.88%  "
# Row 36
$ws.Range("E36").Value = "  -0.18%  "
# Row 37
$ws.Range("D37").Value = "'1.116"
$ws.Range("E37").Value = "  +3.35%  "
# Row 38
$ws.Range("D38").Value = "'0.05212"
$ws.Range("E38").Value = "  +2.84%  "
# Row 39
$ws.Range("E39").Value = "  +1.34%  "
# Row 40
$ws.Range("E40").Value = "  +2.97%  "
# Row 41
$ws.Range("D41").Value = "'0.1642"
$ws.Range("E41").Value = "  +1.84%  "
# Row 42
$ws.Range("D42").Value = "'2.719"
$ws.Range("E42").Value = "  +8.52%  "
# Row 43
$ws.Range("D43").Value = "'6.428"
$ws.Range("E43").Value = "  +7.54%  "
# Row 44
$ws.Range("D44").Value = "'8.265"
$ws.Range("E44").Value = "  +4.19%  "
# Row 45
$ws.Range("D45").Value = "'106.71"
$ws.Range("E45").Value = "  +2.40%  "
# Row 46
$ws.Range("D46").Value = "'10.33"
$ws.Range("E46").Value = "  +1.99%  "
# Row 47
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.77%  "
# Row 48
$ws.Range("D48").Value = "'1.656"
$ws.Range("E48").Value = "  +5.30%  "
# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06288"
$ws.Range("E49").Value = "  +0.28%  "
# Row 50
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "'0.4554"
$ws.Range("E50").Value = "  +1.89%  "
# Row 51
$ws.Range("D51").Value = "'1.800"
$ws.Range("E51").Value = "  +6.46%  "
